$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 119
$ws1.Range("F4").Value = 5128
$ws1.Range("F5").Value = 376
$ws1.Range("F7").Value = 300
$ws1.Range("F8").Value = 778
$ws1.Range("F9").Value = 256

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 119
$ws4.Range("F4").Value = 5128
$ws4.Range("F5").Value = 376
$ws4.Range("F7").Value = 300
$ws4.Range("F8").Value = 778
$ws4.Range("F10").Value = 256
